# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest counts from the regenerated gh-pages data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 64
$ws1.Range("F5").Value  = 1020
$ws1.Range("F8").Value  = 192
$ws1.Range("F9").Value  = 370
$ws1.Range("F10").Value = 3
$ws1.Range("F12").Value = 523
$ws1.Range("F14").Value = 12258
$ws1.Range("F15").Value = 54
$ws1.Range("F16").Value = 5458

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 64
$ws4.Range("F7").Value  = 1020
$ws4.Range("F10").Value = 192
$ws4.Range("F11").Value = 370
$ws4.Range("F12").Value = 3
$ws4.Range("F14").Value = 523
$ws4.Range("F16").Value = 12258
$ws4.Range("F18").Value = 54
$ws4.Range("F19").Value = 5458
